$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (27) down onto
# the two new rows, so the new cells reuse the existing cellXfs entries
# (bold/border/center for column A, the date/time number format for
# column E) instead of Excel fabricating new style records.
$ws.Range("A27").Copy()
$ws.Range("A28:A29").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E28:E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 28 : Manchester 62 0 x 7 Lions Gibraltar ---
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "gibraltar"
$ws.Cells.Item(28, 3).Value = "national-league"
$ws.Cells.Item(28, 4).Value = "2023-2024"
$ws.Cells.Item(28, 5).Value = 45235.6875
$ws.Cells.Item(28, 6).Value = "Manchester 62"
$ws.Cells.Item(28, 7).Value = 7
$ws.Cells.Item(28, 8).Value = "Lions Gibraltar"
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 1.12
$ws.Cells.Item(28, 11).Value = "05/11/2023 11:47"
$ws.Cells.Item(28, 12).Value = 1.15
$ws.Cells.Item(28, 13).Value = "05/11/2023 15:38"
$ws.Cells.Item(28, 14).Value = 8.130000000000001
$ws.Cells.Item(28, 15).Value = "05/11/2023 11:47"
$ws.Cells.Item(28, 16).Value = 8.539999999999999
$ws.Cells.Item(28, 17).Value = "05/11/2023 15:38"
$ws.Cells.Item(28, 18).Value = 9.77
$ws.Cells.Item(28, 19).Value = "05/11/2023 11:47"
$ws.Cells.Item(28, 20).Value = 8.43
$ws.Cells.Item(28, 21).Value = "05/11/2023 15:38"
$ws.Cells.Item(28, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/manchester-62-lions-gibraltar/hpni7Ahr/"

# --- Row 29 : Europa Point 0 x 1 Europa FC ---
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "gibraltar"
$ws.Cells.Item(29, 3).Value = "national-league"
$ws.Cells.Item(29, 4).Value = "2023-2024"
$ws.Cells.Item(29, 5).Value = 45235.8125
$ws.Cells.Item(29, 6).Value = "Europa Point"
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = "Europa FC"
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 3.49
$ws.Cells.Item(29, 11).Value = "05/11/2023 11:47"
$ws.Cells.Item(29, 12).Value = 4.69
$ws.Cells.Item(29, 13).Value = "05/11/2023 19:30"
$ws.Cells.Item(29, 14).Value = 3.94
$ws.Cells.Item(29, 15).Value = "05/11/2023 11:47"
$ws.Cells.Item(29, 16).Value = 4.45
$ws.Cells.Item(29, 17).Value = "05/11/2023 19:30"
$ws.Cells.Item(29, 18).Value = 1.75
$ws.Cells.Item(29, 19).Value = "05/11/2023 11:47"
$ws.Cells.Item(29, 20).Value = 1.5
$ws.Cells.Item(29, 21).Value = "05/11/2023 19:30"
$ws.Cells.Item(29, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/europa-point-europa-fc/Wzoe6U7l/"
